$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows (34-36) replicating the existing data pattern,
# each for regcntr_id 10005 with incrementing usr_id values.
$newRows = @(
    @(10005, 110033),
    @(10005, 110034),
    @(10005, 110035)
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Move the view/selection similarly to what Excel records after scrolling
# to the bottom and selecting the rows below the data.
$ws.Range("A37:XFD1048576").Select()
